# Auto-generated edit script
# Applies new voltage magnitude (vm_pu) values for the 380 kV case
# Column B: slack bus voltage set-point 1.05 -> 1.02
# Columns C-F, I-N: recomputed power-flow results for rows 2-25 (bus indices 0-23)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2 (bus index 0)
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.007741359813101
$ws.Cells.Item(2, 4).Value = 1.010450795242707
$ws.Cells.Item(2, 5).Value = 1.010378232713478
$ws.Cells.Item(2, 6).Value = 1.005934661921953
$ws.Cells.Item(2, 9).Value = 1.023594999628091
$ws.Cells.Item(2, 10).Value = 1.013011236820697
$ws.Cells.Item(2, 11).Value = 1.013320948185163
$ws.Cells.Item(2, 12).Value = 1.013248605009605
$ws.Cells.Item(2, 13).Value = 1.008818532729725
$ws.Cells.Item(2, 14).Value = 1.008424484006033

# Row 3 (bus index 1)
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.008885489017892
$ws.Cells.Item(3, 4).Value = 1.011447612279555
$ws.Cells.Item(3, 5).Value = 1.011354631784088
$ws.Cells.Item(3, 6).Value = 1.007732031210571
$ws.Cells.Item(3, 9).Value = 1.023504579208683
$ws.Cells.Item(3, 10).Value = 1.013785937500267
$ws.Cells.Item(3, 11).Value = 1.014121302574718
$ws.Cells.Item(3, 12).Value = 1.014028582638701
$ws.Cells.Item(3, 13).Value = 1.010416174298673
$ws.Cells.Item(3, 14).Value = 1.008689512514675

# Row 4 (bus index 2)
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.00962571700567
$ws.Cells.Item(4, 4).Value = 1.012092817872978
$ws.Cells.Item(4, 5).Value = 1.01198666391706
$ws.Cells.Item(4, 6).Value = 1.008894589360257
$ws.Cells.Item(4, 9).Value = 1.023443335729026
$ws.Cells.Item(4, 10).Value = 1.014286659102448
$ws.Cells.Item(4, 11).Value = 1.014638783828894
$ws.Cells.Item(4, 12).Value = 1.014532912232917
$ws.Cells.Item(4, 13).Value = 1.011449090507858
$ws.Cells.Item(4, 14).Value = 1.008860573501175

# Row 5 (bus index 3)
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.009936887534993
$ws.Cells.Item(5, 4).Value = 1.012364111420605
$ws.Cells.Item(5, 5).Value = 1.012252428344766
$ws.Cells.Item(5, 6).Value = 1.00938322889415
$ws.Cells.Item(5, 9).Value = 1.023416932628374
$ws.Cells.Item(5, 10).Value = 1.014497029797006
$ws.Cells.Item(5, 11).Value = 1.014856238440872
$ws.Cells.Item(5, 12).Value = 1.014744845756573
$ws.Cells.Item(5, 13).Value = 1.011883130412384
$ws.Cells.Item(5, 14).Value = 1.008932384960895

# Row 6 (bus index 4)
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.009989133262976
$ws.Cells.Item(6, 4).Value = 1.012409665728678
$ws.Cells.Item(6, 5).Value = 1.012297054812176
$ws.Cells.Item(6, 6).Value = 1.009465268058942
$ws.Cells.Item(6, 9).Value = 1.023412460900047
$ws.Cells.Item(6, 10).Value = 1.014532344233848
$ws.Cells.Item(6, 11).Value = 1.014892744532218
$ws.Cells.Item(6, 12).Value = 1.014780425290233
$ws.Cells.Item(6, 13).Value = 1.011955996266641
$ws.Cells.Item(6, 14).Value = 1.008944436417451

# Row 7 (bus index 5)
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.009629874961979
$ws.Cells.Item(7, 4).Value = 1.012096442715623
$ws.Cells.Item(7, 5).Value = 1.011990214844175
$ws.Cells.Item(7, 6).Value = 1.008901118965581
$ws.Cells.Item(7, 9).Value = 1.023442985510381
$ws.Cells.Item(7, 10).Value = 1.014289470605075
$ws.Cells.Item(7, 11).Value = 1.014641689837406
$ws.Cells.Item(7, 12).Value = 1.014535744437633
$ws.Cells.Item(7, 13).Value = 1.011454890937522
$ws.Cells.Item(7, 14).Value = 1.008861533451982

# Row 8 (bus index 6)
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.008128045961338
$ws.Cells.Item(8, 4).Value = 1.010787634154178
$ws.Cells.Item(8, 5).Value = 1.010708163237113
$ws.Cells.Item(8, 6).Value = 1.006542192680029
$ws.Cells.Item(8, 9).Value = 1.023565007130014
$ws.Cells.Item(8, 10).Value = 1.013273167878998
$ws.Cells.Item(8, 11).Value = 1.013591515823637
$ws.Cells.Item(8, 12).Value = 1.013512279211208
$ws.Cells.Item(8, 13).Value = 1.009358647616567
$ws.Cells.Item(8, 14).Value = 1.008514140949303

# Row 9 (bus index 7)
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.005480745576251
$ws.Cells.Item(9, 4).Value = 1.008482791098277
$ws.Cells.Item(9, 5).Value = 1.00845077064082
$ws.Cells.Item(9, 6).Value = 1.00238148300482
$ws.Cells.Item(9, 9).Value = 1.023759139049229
$ws.Cells.Item(9, 10).Value = 1.01147793231855
$ws.Cells.Item(9, 11).Value = 1.011737826317246
$ws.Cells.Item(9, 12).Value = 1.011705916681197
$ws.Cells.Item(9, 13).Value = 1.005657775411462
$ws.Cells.Item(9, 14).Value = 1.007898677791794

# Row 10 (bus index 8)
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.003715114672957
$ws.Cells.Item(10, 4).Value = 1.006947094622776
$ws.Cells.Item(10, 5).Value = 1.006946921126703
$ws.Cells.Item(10, 6).Value = 0.9996043609359114
$ws.Cells.Item(10, 9).Value = 1.023874603872166
$ws.Cells.Item(10, 10).Value = 1.01027805046223
$ws.Cells.Item(10, 11).Value = 1.010499811863725
$ws.Cells.Item(10, 12).Value = 1.010499639027784
$ws.Cells.Item(10, 13).Value = 1.003185244096014
$ws.Cells.Item(10, 14).Value = 1.007486114562519

# Row 11 (bus index 9)
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.002950357702456
$ws.Cells.Item(11, 4).Value = 1.006282302090868
$ws.Cells.Item(11, 5).Value = 1.006295971652534
$ws.Cells.Item(11, 6).Value = 0.998400907678803
$ws.Cells.Item(11, 9).Value = 1.023921309042974
$ws.Cells.Item(11, 10).Value = 1.009757737288748
$ws.Cells.Item(11, 11).Value = 1.009963187190324
$ws.Cells.Item(11, 12).Value = 1.009976802697977
$ws.Cells.Item(11, 13).Value = 1.002113232761298
$ws.Cells.Item(11, 14).Value = 1.007306929117971

# Row 12 (bus index 10)
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.002666255427317
$ws.Cells.Item(12, 4).Value = 1.006035392536111
$ws.Cells.Item(12, 5).Value = 1.006054212209226
$ws.Cells.Item(12, 6).Value = 0.9979537378880736
$ws.Cells.Item(12, 9).Value = 1.023938164268129
$ws.Cells.Item(12, 10).Value = 1.009564353877718
$ws.Cells.Item(12, 11).Value = 1.0097637751944
$ws.Cells.Item(12, 12).Value = 1.009782519401926
$ws.Cells.Item(12, 13).Value = 1.001714820554971
$ws.Cells.Item(12, 14).Value = 1.007240289530172

# Row 13 (bus index 11)
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.002727198080188
$ws.Cells.Item(13, 4).Value = 1.006088354423225
$ws.Cells.Item(13, 5).Value = 1.006106069006355
$ws.Cells.Item(13, 6).Value = 0.9980496644779331
$ws.Cells.Item(13, 9).Value = 1.023934571058996
$ws.Cells.Item(13, 10).Value = 1.009605840558904
$ws.Cells.Item(13, 11).Value = 1.009806553656861
$ws.Cells.Item(13, 12).Value = 1.009824197417062
$ws.Cells.Item(13, 13).Value = 1.0018002914309
$ws.Cells.Item(13, 14).Value = 1.007254587674138

# Row 14 (bus index 12)
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.002926874503961
$ws.Cells.Item(14, 4).Value = 1.006261891993281
$ws.Cells.Item(14, 5).Value = 1.006275987096645
$ws.Cells.Item(14, 6).Value = 0.9983639476905006
$ws.Cells.Item(14, 9).Value = 1.023922712353276
$ws.Cells.Item(14, 10).Value = 1.009741754533931
$ws.Cells.Item(14, 11).Value = 1.009946705493141
$ws.Cells.Item(14, 12).Value = 1.009960744795007
$ws.Cells.Item(14, 13).Value = 1.002080304423726
$ws.Cells.Item(14, 14).Value = 1.007301422347768

# Row 15 (bus index 13)
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.003049896670781
$ws.Cells.Item(15, 4).Value = 1.006368817331356
$ws.Cells.Item(15, 5).Value = 1.006380683444271
$ws.Cells.Item(15, 6).Value = 0.998557567219246
$ws.Cells.Item(15, 9).Value = 1.023915340507297
$ws.Cells.Item(15, 10).Value = 1.00982548018888
$ws.Cells.Item(15, 11).Value = 1.010033046245149
$ws.Cells.Item(15, 12).Value = 1.010044865666519
$ws.Cells.Item(15, 13).Value = 1.002252800264956
$ws.Cells.Item(15, 14).Value = 1.007330267834747

# Row 16 (bus index 14)
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.003765864002442
$ws.Cells.Item(16, 4).Value = 1.006991218197981
$ws.Cells.Item(16, 5).Value = 1.006990127120171
$ws.Cells.Item(16, 6).Value = 0.9996842092912569
$ws.Cells.Item(16, 9).Value = 1.023871434968932
$ws.Cells.Item(16, 10).Value = 1.010312565802187
$ws.Cells.Item(16, 11).Value = 1.010535413975112
$ws.Cells.Item(16, 12).Value = 1.010534327037911
$ws.Cells.Item(16, 13).Value = 1.003256359774522
$ws.Cells.Item(16, 14).Value = 1.00749799502394

# Row 17 (bus index 15)
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.004214908215161
$ws.Cells.Item(17, 4).Value = 1.007381678918184
$ws.Cells.Item(17, 5).Value = 1.007372474500553
$ws.Cells.Item(17, 6).Value = 1.00039066149902
$ws.Cells.Item(17, 9).Value = 1.023843013862359
$ws.Cells.Item(17, 10).Value = 1.01061789761435
$ws.Cells.Item(17, 11).Value = 1.01085038555817
$ws.Cells.Item(17, 12).Value = 1.010841215264662
$ws.Cells.Item(17, 13).Value = 1.003885486680965
$ws.Cells.Item(17, 14).Value = 1.007603060116528

# Row 18 (bus index 16)
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.004476806420504
$ws.Cells.Item(18, 4).Value = 1.007609445004532
$ws.Cells.Item(18, 5).Value = 1.007595513323928
$ws.Cells.Item(18, 6).Value = 1.000802633153154
$ws.Cells.Item(18, 9).Value = 1.023826118601784
$ws.Cells.Item(18, 10).Value = 1.010795919667981
$ws.Cells.Item(18, 11).Value = 1.011034049629468
$ws.Cells.Item(18, 12).Value = 1.011020168887485
$ws.Cells.Item(18, 13).Value = 1.004252312479569
$ws.Cells.Item(18, 14).Value = 1.007664290470718

# Row 19 (bus index 17)
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.004566103363046
$ws.Cells.Item(19, 4).Value = 1.007687110294921
$ws.Cells.Item(19, 5).Value = 1.007671567600544
$ws.Cells.Item(19, 6).Value = 1.000943089951269
$ws.Cells.Item(19, 9).Value = 1.023820303822167
$ws.Cells.Item(19, 10).Value = 1.01085660830383
$ws.Cells.Item(19, 11).Value = 1.011096665316714
$ws.Cells.Item(19, 12).Value = 1.011081179180929
$ws.Cells.Item(19, 13).Value = 1.004377368315579
$ws.Cells.Item(19, 14).Value = 1.007685159579239

# Row 20 (bus index 18)
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.004166732289501
$ws.Cells.Item(20, 4).Value = 1.007339784441221
$ws.Cells.Item(20, 5).Value = 1.007331449978018
$ws.Cells.Item(20, 6).Value = 1.000314875259818
$ws.Cells.Item(20, 9).Value = 1.023846096030192
$ws.Cells.Item(20, 10).Value = 1.010585145950745
$ws.Cells.Item(20, 11).Value = 1.01081659764068
$ws.Cells.Item(20, 12).Value = 1.010808294154659
$ws.Cells.Item(20, 13).Value = 1.003818001184111
$ws.Cells.Item(20, 14).Value = 1.007591793044342

# Row 21 (bus index 19)
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.002868075835822
$ws.Cells.Item(21, 4).Value = 1.006210788859491
$ws.Cells.Item(21, 5).Value = 1.00622594958189
$ws.Cells.Item(21, 6).Value = 0.9982714034586954
$ws.Cells.Item(21, 9).Value = 1.02392621804881
$ws.Cells.Item(21, 10).Value = 1.009701734474589
$ws.Cells.Item(21, 11).Value = 1.009905436653439
$ws.Cells.Item(21, 12).Value = 1.009920537181223
$ws.Cells.Item(21, 13).Value = 1.001997853721787
$ws.Cells.Item(21, 14).Value = 1.007287632975731

# Row 22 (bus index 20)
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.002051340423527
$ws.Cells.Item(22, 4).Value = 1.005501082244604
$ws.Cells.Item(22, 5).Value = 1.00553106277097
$ws.Cells.Item(22, 6).Value = 0.9969856964212369
$ws.Cells.Item(22, 9).Value = 1.023973741382756
$ws.Cells.Item(22, 10).Value = 1.009145626728355
$ws.Cells.Item(22, 11).Value = 1.0093320566831
$ws.Cells.Item(22, 12).Value = 1.009361913377761
$ws.Cells.Item(22, 13).Value = 1.0008521796193
$ws.Cells.Item(22, 14).Value = 1.007095919966662

# Row 23 (bus index 21)
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.002484329161242
$ws.Cells.Item(23, 4).Value = 1.005877298825826
$ws.Cells.Item(23, 5).Value = 1.00589941852797
$ws.Cells.Item(23, 6).Value = 0.9976673630241029
$ws.Cells.Item(23, 9).Value = 1.023948818288664
$ws.Cells.Item(23, 10).Value = 1.009440494421733
$ws.Cells.Item(23, 11).Value = 1.009636064068379
$ws.Cells.Item(23, 12).Value = 1.009658094284072
$ws.Cells.Item(23, 13).Value = 1.001459647627042
$ws.Cells.Item(23, 14).Value = 1.007197595927036

# Row 24 (bus index 22)
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.004188500983481
$ws.Cells.Item(24, 4).Value = 1.00735871469962
$ws.Cells.Item(24, 5).Value = 1.007349987125685
$ws.Cells.Item(24, 6).Value = 1.000349120075986
$ws.Cells.Item(24, 9).Value = 1.023844704313072
$ws.Cells.Item(24, 10).Value = 1.010599945242685
$ws.Cells.Item(24, 11).Value = 1.010831865111529
$ws.Cells.Item(24, 12).Value = 1.010823169938823
$ws.Cells.Item(24, 13).Value = 1.003848495386098
$ws.Cells.Item(24, 14).Value = 1.00759688431099

# Row 25 (bus index 23)
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.006165258446153
$ws.Cells.Item(25, 4).Value = 1.009078487903274
$ws.Cells.Item(25, 5).Value = 1.009034163554825
$ws.Cells.Item(25, 6).Value = 1.003457663494368
$ws.Cells.Item(25, 9).Value = 1.02371141742819
$ws.Cells.Item(25, 10).Value = 1.011942574428572
$ws.Cells.Item(25, 11).Value = 1.01221743333349
$ws.Cells.Item(25, 12).Value = 1.012173256498776
$ws.Cells.Item(25, 13).Value = 1.006615429486483
$ws.Cells.Item(25, 14).Value = 1.008058185146165

Write-Output "Applied 380 kV case vm_pu updates"
